$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I10").Value = 3.4
$ws.Range("J10").Value = 3.4
$ws.Range("K10").Value = 1.83
$ws.Range("M10").Value = 1.13
$ws.Range("N10").Value = 6
$ws.Range("O10").Value = 1.57
$ws.Range("P10").Value = 2.25
$ws.Range("Q10").Value = 2.88
$ws.Range("R10").Value = 1.4

$ws.Range("AA10").Value = 26
$ws.Range("AF10").Value = 81
$ws.Range("AJ10").Value = 34
$ws.Range("AP10").Value = 34
$ws.Range("AR10").Value = 101
